$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with refreshed scores ---
$ws.Range("F94").Value = 9.476147036988191
$ws.Range("I94").Value = 6.960549220560061
$ws.Range("P94").Value = 33.2355770570365
$ws.Range("Q94").Value = 36.23999171389823

$ws.Range("C98").Value = 8.730719641740396
$ws.Range("E98").Value = 10
$ws.Range("F98").Value = 9.426184165891407
$ws.Range("G98").Value = 9.893083761673136
$ws.Range("H98").Value = 10
$ws.Range("I98").Value = 9.820660211092806
$ws.Range("J98").Value = 9.983267724329723
$ws.Range("L98").Value = 9.077626808507636
$ws.Range("M98").Value = 8.123058660957394
$ws.Range("P98").Value = 46.56752227546373
$ws.Range("Q98").Value = 38.48707869872877

$ws.Range("C99").Value = 6.679803770256861
$ws.Range("E99").Value = 7.281635567666262
$ws.Range("F99").Value = 5.866683268186358
$ws.Range("H99").Value = 10
$ws.Range("I99").Value = 8.90625
$ws.Range("J99").Value = 7.280332666343885
$ws.Range("M99").Value = 5.647878198898607
$ws.Range("P99").Value = 40.40755240094282
$ws.Range("Q99").Value = 33.14701593453024

$ws.Range("F100").Value = 7.466666666666665
$ws.Range("I100").Value = 8.166666666666666
$ws.Range("P100").Value = 44.43333333333333
$ws.Range("Q100").Value = 30.16666666666666

# --- Append new rows 102-105 for 2025-02-26 ---
$dateRange = $ws.Range("A102:A105")
$dateRange.NumberFormat = "@"
$ws.Range("A102").Value = "2025-02-26"
$ws.Range("B102").Value = "abs_activity"
$ws.Range("C102").Value = 7.876614166839925
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 7.012697428973726
$ws.Range("F102").Value = 8.884237811123356
$ws.Range("G102").Value = 8.70407878200314
$ws.Range("H102").Value = 10
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 8.638168564234581
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 5.206230415513048
$ws.Range("M102").Value = 9.15910428370217
$ws.Range("N102").Value = 0
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 32.75249466151897
$ws.Range("Q102").Value = 32.72863679087099

$ws.Range("A103").Value = "2025-02-26"
$ws.Range("B103").Value = "rel_activity"
$ws.Range("C103").Value = 5.179601660970858
$ws.Range("D103").Value = 5
$ws.Range("E103").Value = 0
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 5.085999312005504
$ws.Range("H103").Value = 10
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 5.916666666666667
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 5.615569522091636
$ws.Range("N103").Value = 5
$ws.Range("O103").Value = 5
$ws.Range("P103").Value = 20.881170495068
$ws.Range("Q103").Value = 25.91666666666667

$ws.Range("A104").Value = "2025-02-26"
$ws.Range("B104").Value = "abs_sleep"
$ws.Range("C104").Value = 8.9
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = 4.866666666666667
$ws.Range("G104").Value = 9.366666666666667
$ws.Range("H104").Value = 4.399999999999999
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 10
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 10
$ws.Range("M104").Value = 9.7
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 0
$ws.Range("P104").Value = 37.96666666666667
$ws.Range("Q104").Value = 29.26666666666667

$ws.Range("A105").Value = "2025-02-26"
$ws.Range("B105").Value = "rel_sleep"
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 0
$ws.Range("E105").Value = 8.60717474553067
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 8.155773955773958
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 10
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("O105").Value = 0
$ws.Range("P105").Value = 8.60717474553067
$ws.Range("Q105").Value = 18.15577395577396

$dateRange.ClearFormats()
